$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns AP1:AS1
$ws.Range("AP1").Value = "ACCEPTED"
$ws.Range("AQ1").Value = "PAID"
$ws.Range("AR1").Value = "HOLD"
$ws.Range("AS1").Value = "REJECTED"

# Copy the header style (s="1") from AO1 to the new header cells
$ws.Range("AO1").Copy()
$ws.Range("AP1:AS1").PasteSpecial(-4122)  # xlPasteFormats

# Set AI/AJ to "HOLD" for rows 5 through 16
for ($r = 5; $r -le 16; $r++) {
    $ws.Cells.Item($r, 35).Value = "HOLD"  # AI
    $ws.Cells.Item($r, 36).Value = "HOLD"  # AJ
}

# Fill boolean values for AP:AS for rows 2 through 16
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 42).Value = $false  # AP
    $ws.Cells.Item($r, 43).Value = $false  # AQ
    if ($r -ge 5) {
        $ws.Cells.Item($r, 44).Value = $true   # AR
    } else {
        $ws.Cells.Item($r, 44).Value = $false  # AR
    }
    $ws.Cells.Item($r, 45).Value = $false  # AS
}
